# Translate table/column names from German to English in the "Admin" etc.
# database documentation sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A-column section headers and row labels: German -> English translations
$ws.Range("A4").Value = "Id_Admin"
$ws.Range("A5").Value = "User_Name"
$ws.Range("A6").Value = "Password"

$ws.Range("A8").Value = "Trainee"
$ws.Range("A9").Value = "Id_Trainee"
$ws.Range("A10").Value = "Login_Code"
$ws.Range("A11").Value = "Id_Godfather"

$ws.Range("A13").Value = "Godfather"
$ws.Range("A14").Value = "Id_Godfather"
$ws.Range("A15").Value = "Email"
$ws.Range("A16").Value = "Password"
$ws.Range("A17").Value = "Last_Name"
$ws.Range("A18").Value = "First_Name"
$ws.Range("A19").Value = "Id_Location"
$ws.Range("A20").Value = "Description"
$ws.Range("A21").Value = "Picture"
$ws.Range("A22").Value = "Id_Teaching_Type"
$ws.Range("A23").Value = "Id_Job"
$ws.Range("A24").Value = "Hiring_Date"
$ws.Range("B24").Value = "Date"
$ws.Range("C24").ClearContents()
$ws.Range("A25").Value = "Birthday"
$ws.Range("B25").Value = "Date"
$ws.Range("C25").ClearContents()
$ws.Range("A26").Value = "Email"
$ws.Range("A27").Value = "Pick_Text"

$ws.Range("A29").Value = "Location"
$ws.Range("A30").Value = "Id_Location"
$ws.Range("A31").Value = "Location"

$ws.Range("A33").Value = "Teaching_Type"
$ws.Range("A34").Value = "Id_Teaching_Type"
$ws.Range("A35").Value = "Teaching_Type"

$ws.Range("A37").Value = "Job"
$ws.Range("A38").Value = "Id_Job"
$ws.Range("A39").Value = "Job"

$ws.Range("A9").Select()
